$wb = $excel.ActiveWorkbook

# Sheet "展览" (first worksheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 14579
$ws1.Range("F5").Value = 17585
$ws1.Range("F7").Value = 69
$ws1.Range("F8").Value = 54
$ws1.Range("F16").Value = 38
$ws1.Range("F17").Value = 140
$ws1.Range("F19").Value = 1327
$ws1.Range("F23").Value = 211
$ws1.Range("F24").Value = 7224
$ws1.Range("F25").Value = 979
$ws1.Range("F28").Value = 1170
$ws1.Range("F30").Value = 5859
$ws1.Range("F31").Value = 64
$ws1.Range("F33").Value = 135
$ws1.Range("F36").Value = 5078

# Sheet "全部类型" (fourth worksheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 14579
$ws4.Range("F5").Value = 17585
$ws4.Range("F7").Value = 69
$ws4.Range("F8").Value = 54
$ws4.Range("F16").Value = 38
$ws4.Range("F17").Value = 140
$ws4.Range("F19").Value = 1327
$ws4.Range("F24").Value = 211
$ws4.Range("F25").Value = 7224
$ws4.Range("F26").Value = 979
$ws4.Range("F29").Value = 1170
$ws4.Range("F32").Value = 5859
$ws4.Range("F33").Value = 64
$ws4.Range("F35").Value = 135
$ws4.Range("F38").Value = 5078
